$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("118500", "17706585", "6004"),
    @("118518", "17706586", "6020"),
    @("118498", "17706587", "1001"),
    @("118452", "17706588", "1001"),
    @("118518", "17706589", "6012"),
    @("118448", "17706590", "1047"),
    @("118518", "17706591", "1035"),
    @("118452", "17706592", "1150"),
    @("118463", "17707507", "1010"),
    @("118463", "17707512", "1010"),
    @("118463", "17707515", "1010"),
    @("118463", "17707516", "1003"),
    @("118463", "17707517", "1007"),
    @("118464", "17707519", "1010"),
    @("118465", "17707520", "1011")
)

$startRow = 128
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
